$d = $word.ActiveDocument

function New-PkgXml($bodyInner) {
    return '<?xml version="1.0"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' + $bodyInner + '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark from its original location (near "Google...")
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2) Give the blank paragraph (before "Retrieving current location of
#    device:") the ListParagraph style, and move the _GoBack bookmark here.
# ---------------------------------------------------------------------------
$pBlank = $d.Paragraphs.Item(9)
$blankXml = New-PkgXml('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
$null = $pBlank.Range.InsertXML($blankXml)

# ---------------------------------------------------------------------------
# 3) Merge the runs in the "If your app does navigation..." paragraph.
# ---------------------------------------------------------------------------
$rpr = '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr>'

$runA = '<w:r>' + $rpr + "<w:t>If your app does navigation or tracking, you probably want to get the user's location at regular intervals.</w:t></w:r>"
$runB = '<w:r>' + $rpr + '<w:t xml:space="preserve"> You can do this by requesting periodic updates from Location Services.</w:t></w:r>'
$runC = '<w:r>' + $rpr + '<w:t xml:space="preserve"> In response, Location Services automatically updates your app with the best available location, based on the currently-available location providers such as WiFi and GPS.</w:t></w:r>'

$findRng = $d.Content
$findRng.Find.ClearFormatting()
$null = $findRng.Find.Execute("If your app does navigation or tracking*such as WiFi and GPS.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
$targetRng = $d.Range($findRng.Start, $findRng.End)
$null = $targetRng.InsertXML(New-PkgXml('<w:p>' + $runA + $runB + $runC + '</w:p><w:p/>'))

# ---------------------------------------------------------------------------
# 4) Split the "Reverse coding can be done using getFromLocation()..." run.
# ---------------------------------------------------------------------------
$reverseXml = (
  '<w:r>' + $rpr + '<w:t xml:space="preserve">Reverse coding can be done using </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r>' + $rpr + '<w:t>getFromLocation</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rpr + '<w:t>(</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r>' + $rpr + '<w:t>) method of Geocoder class.</w:t></w:r>'
)
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$null = $findRng.Find.Execute("Reverse coding can be done using getFromLocation`(`) method of Geocoder class.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetRng = $d.Range($findRng.Start, $findRng.End)
$null = $targetRng.InsertXML(New-PkgXml('<w:p>' + $reverseXml + '</w:p><w:p/>'))

# ---------------------------------------------------------------------------
# 5) Split the geofence-expiration sentence, adding gramStart/gramEnd marks.
# ---------------------------------------------------------------------------
$geoXml = (
  '<w:r>' + $rpr + '<w:t xml:space="preserve">Location Services treats a geofences as an area rather than as a points and proximity. This allows it to detect when the user enters or exits a geofence. For each geofence, you can ask Location Services to send you entrance events or exit events or both. You can also limit the duration of a geofence by specifying </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r>' + $rpr + '<w:t>an expiration</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r>' + $rpr + '<w:t xml:space="preserve"> duration in milliseconds. After the geofence expires, Location Services automatically removes it. </w:t></w:r>'
)
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$null = $findRng.Find.Execute("Location Services treats a geofences*automatically removes it. ", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
$targetRng = $d.Range($findRng.Start, $findRng.End)
$null = $targetRng.InsertXML(New-PkgXml('<w:p>' + $geoXml + '</w:p><w:p/>'))

# ---------------------------------------------------------------------------
# 6) Split the activity-recognition sentence, marking "geofencing" as a
#    spelling flag.
# ---------------------------------------------------------------------------
$actXml = (
  '<w:r>' + $rpr + "<w:t xml:space=`"preserve`">Activity recognition tries to detect the user's current physical activity, such as walking, driving, or standing still. Requests for updates go through an activity recognition client, which, while different from the location client used by location or </w:t></w:r>" +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rpr + '<w:t>geofencing</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rpr + '<w:t>, follows a similar pattern. Based on the update interval you choose, Location Services sends out activity information containing one or more possible activities and the confidence level for each one.</w:t></w:r>'
)
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$null = $findRng.Find.Execute("Activity recognition tries to detect*confidence level for each one.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
$targetRng = $d.Range($findRng.Start, $findRng.End)
$null = $targetRng.InsertXML(New-PkgXml('<w:p>' + $actXml + '</w:p><w:p/>'))

Write-Host "All edits applied."
